$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "El objetivo de " + "Rush " -> single run "El objetivo de Rush "
#    These two runs share identical rPr and are a self-contained
#    contiguous pair (bounded on both sides by <w:proofErr/> markers), so
#    a same-text Find/Replace coalesces them into one run without
#    touching anything else.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("El objetivo de Rush ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "El objetivo de Rush ", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) " Paisa" + " es intentar representar un poco de l" -> single run,
#    while the run that follows ("as situaciones...") must stay a
#    separate run. A Find/Replace across this boundary ends up
#    coalescing that following run too (engine quirk), so the merge is
#    done surgically: locate the exact text range, then splice in the
#    desired two-run OOXML via Range.InsertXML.
# ---------------------------------------------------------------------
$rFind = $d.Content
$rFind.Find.Execute(" Paisa es intentar representar un poco de l", $true, $false, $false, $false, $false, `
                     $true, 1, $false, $null, 0) | Out-Null
$matchStart = $rFind.Start
$para = $rFind.Paragraphs(1).Range
$paraTextEnd = $para.End - 1
$spliceRange = $d.Range($matchStart, $paraTextEnd)
$spliceRange.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve"> Paisa es intentar representar un poco de l</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">as situaciones que viven las personas que utilizan constantemente el sistema metro en la ciudad llevándolo a un plano mas surrealista como lo puede llegar a ser un videojuego, la idea es que el jugador mantenga un balance entre distintas “barras de vida” que disminuirán o aumentaran en función a las decisiones que tome el jugador durante el viaje.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# ---------------------------------------------------------------------
# 3) Remove the trailing paragraph "El motivo por el cual escojo este
#    evento..." entirely (including its paragraph mark).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("El motivo por el cual escojo este evento")) {
        $p.Range.Delete()
        break
    }
}
